# "Fixed issue with cloned metadata"
#
# Inserts a new "Projectiles" worksheet right after "Weapons" (so the tab
# order becomes Weapons, Projectiles, Potion, Enemies, Armor), fills it with
# the arrow items, makes it the active sheet/selection, and updates the
# selection left behind on the Weapons sheet.

$wb = $excel.ActiveWorkbook

$weapons = $wb.Worksheets.Item("Weapons")

# Insert the new sheet immediately after "Weapons".
$projectiles = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $weapons)
$projectiles.Name = "Projectiles"

# Header row (same column layout/order as the other item sheets).
$projectiles.Range("A1").Value = "Name"
$projectiles.Range("B1").Value = "Level"
$projectiles.Range("C1").Value = "Damage"
$projectiles.Range("D1").Value = "Weight"
$projectiles.Range("E1").Value = "Bulk "
$projectiles.Range("F1").Value = "Value"
$projectiles.Range("A1:F1").Font.Bold = $true

# Data rows.
$projectiles.Range("A2").Value = "Wooden Arrow"
$projectiles.Range("B2").Value = 1
$projectiles.Range("C2").Value = "1d4"
$projectiles.Range("D2").Value = 1
$projectiles.Range("E2").Value = 2
$projectiles.Range("F2").Value = 1

$projectiles.Range("A3").Value = "Steel Arrow"
$projectiles.Range("B3").Value = 1
$projectiles.Range("C3").Value = "1d5"
$projectiles.Range("D3").Value = 1
$projectiles.Range("E3").Value = 2
$projectiles.Range("F3").Value = 2

$projectiles.Range("A4").Value = "Silver Arrow"
$projectiles.Range("B4").Value = 1
$projectiles.Range("C4").Value = "2d3"
$projectiles.Range("D4").Value = 1
$projectiles.Range("E4").Value = 2
$projectiles.Range("F4").Value = 4

$projectiles.Range("A5").Value = "Flaming Arrow"
$projectiles.Range("B5").Value = 2
$projectiles.Range("C5").Value = "2d6"
$projectiles.Range("D5").Value = 1
$projectiles.Range("E5").Value = 1
$projectiles.Range("F5").Value = 10

$projectiles.Columns("A:F").AutoFit() | Out-Null

# Leave the Weapons sheet's selection where the author left it.
$weapons.Range("E11:F11").Select() | Out-Null

# Projectiles ends up the active tab, selection on B5.
$projectiles.Activate() | Out-Null
$projectiles.Range("B5").Select() | Out-Null
